$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 45654
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(3, 2).Value = 45654
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(4, 2).Value = 45641
$ws.Cells.Item(4, 3).Value = 0.9997152494852587
$ws.Cells.Item(5, 2).Value = 45654
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(6, 2).Value = 45654
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(7, 2).Value = 45595
$ws.Cells.Item(7, 3).Value = 0.9987076707407894
$ws.Cells.Item(8, 2).Value = 45654
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(9, 2).Value = 297
$ws.Cells.Item(9, 3).Value = 0.006505454067551584
$ws.Cells.Item(10, 2).Value = 45654
$ws.Cells.Item(10, 3).Value = 1
$ws.Cells.Item(11, 2).Value = 45654
$ws.Cells.Item(11, 3).Value = 1
$ws.Cells.Item(12, 2).Value = 45654
$ws.Cells.Item(12, 3).Value = 1
$ws.Cells.Item(13, 2).Value = 45654
$ws.Cells.Item(13, 3).Value = 1
$ws.Cells.Item(14, 2).Value = 45654
$ws.Cells.Item(14, 3).Value = 1
$ws.Cells.Item(15, 2).Value = 45654
$ws.Cells.Item(15, 3).Value = 1
$ws.Cells.Item(16, 2).Value = 45654
$ws.Cells.Item(16, 3).Value = 1
$ws.Cells.Item(17, 2).Value = 45654
$ws.Cells.Item(17, 3).Value = 1
$ws.Cells.Item(18, 2).Value = 1281
$ws.Cells.Item(18, 3).Value = 0.0280588776448942
$ws.Cells.Item(19, 2).Value = 45603
$ws.Cells.Item(19, 3).Value = 0.9988829018267841
$ws.Cells.Item(20, 2).Value = 2911
$ws.Cells.Item(20, 3).Value = 0.06376221141630525
$ws.Cells.Item(21, 2).Value = 672
$ws.Cells.Item(21, 3).Value = 0.01471941122355106
$ws.Cells.Item(22, 2).Value = 7210
$ws.Cells.Item(22, 3).Value = 0.1579270162526832
$ws.Cells.Item(23, 2).Value = 5237
$ws.Cells.Item(23, 3).Value = 0.1147106496692513
$ws.Cells.Item(24, 2).Value = 1016
$ws.Cells.Item(24, 3).Value = 0.02225434792132124
$ws.Cells.Item(25, 2).Value = 32969
$ws.Cells.Item(25, 3).Value = 0.7221492092697245
$ws.Cells.Item(26, 2).Value = 702
$ws.Cells.Item(26, 3).Value = 0.01537652779603102
$ws.Cells.Item(28, 2).Value = 45621
$ws.Cells.Item(28, 3).Value = 0.9992771717702721
$ws.Cells.Item(29, 2).Value = 45633
$ws.Cells.Item(29, 3).Value = 0.999540018399264
$ws.Cells.Item(30, 2).Value = 45624
$ws.Cells.Item(30, 3).Value = 0.9993428834275201
$ws.Cells.Item(32, 2).Value = 45654
$ws.Cells.Item(32, 3).Value = 1
$ws.Cells.Item(33, 2).Value = 36021
$ws.Cells.Item(33, 3).Value = 0.7889998685766855
$ws.Cells.Item(34, 2).Value = 44458
$ws.Cells.Item(34, 3).Value = 0.973802952643799
$ws.Cells.Item(35, 2).Value = 45595
$ws.Cells.Item(35, 3).Value = 0.9987076707407894
$ws.Cells.Item(36, 2).Value = 45268
$ws.Cells.Item(36, 3).Value = 0.9915451001007579
$ws.Cells.Item(37, 2).Value = 45537
$ws.Cells.Item(37, 3).Value = 0.9974372453673281
$ws.Cells.Item(38, 2).Value = 45372
$ws.Cells.Item(38, 3).Value = 0.9938231042186884
$ws.Cells.Item(39, 2).Value = 45642
$ws.Cells.Item(39, 3).Value = 0.999737153371008
$ws.Cells.Item(40, 2).Value = 45624
$ws.Cells.Item(40, 3).Value = 0.9993428834275201
$ws.Cells.Item(41, 2).Value = 45619
$ws.Cells.Item(41, 3).Value = 0.9992333639987734
$ws.Cells.Item(42, 2).Value = 1221
$ws.Cells.Item(42, 3).Value = 0.02674464449993429
$ws.Cells.Item(43, 2).Value = 78
$ws.Cells.Item(43, 3).Value = 0.001708503088447891
$ws.Cells.Item(44, 2).Value = 843
$ws.Cells.Item(44, 3).Value = 0.01846497568668682
$ws.Cells.Item(45, 2).Value = 121
$ws.Cells.Item(45, 3).Value = 0.002650370175669164
$ws.Cells.Item(46, 2).Value = 300
$ws.Cells.Item(46, 3).Value = 0.006571165724799579
$ws.Cells.Item(47, 2).Value = 575
$ws.Cells.Item(47, 3).Value = 0.01259473430586586
$ws.Cells.Item(48, 2).Value = 1358
$ws.Cells.Item(48, 3).Value = 0.02974547684759276
$ws.Cells.Item(49, 2).Value = 83
$ws.Cells.Item(49, 3).Value = 0.00181802251719455
$ws.Cells.Item(50, 2).Value = 300
$ws.Cells.Item(50, 3).Value = 0.006571165724799579
$ws.Cells.Item(51, 2).Value = 288
$ws.Cells.Item(51, 3).Value = 0.006308319095807596
$ws.Cells.Item(52, 2).Value = 3384
$ws.Cells.Item(52, 3).Value = 0.07412274937573926
$ws.Cells.Item(53, 2).Value = 45534
$ws.Cells.Item(53, 3).Value = 0.9973715337100801
$ws.Cells.Item(54, 2).Value = 45608
$ws.Cells.Item(54, 3).Value = 0.9989924212555308
$ws.Cells.Item(55, 2).Value = 45619
$ws.Cells.Item(55, 3).Value = 0.9992333639987734
$ws.Cells.Item(56, 2).Value = 45607
$ws.Cells.Item(56, 3).Value = 0.9989705173697814
$ws.Cells.Item(57, 2).Value = 45607
$ws.Cells.Item(57, 3).Value = 0.9989705173697814
$ws.Cells.Item(58, 2).Value = 45598
$ws.Cells.Item(58, 3).Value = 0.9987733823980374
$ws.Cells.Item(59, 2).Value = 45605
$ws.Cells.Item(59, 3).Value = 0.9989267095982828
$ws.Cells.Item(60, 2).Value = 45610
$ws.Cells.Item(60, 3).Value = 0.9990362290270294
$ws.Cells.Item(61, 2).Value = 45343
$ws.Cells.Item(61, 3).Value = 0.9931878915319577
$ws.Cells.Item(62, 2).Value = 36325
$ws.Cells.Item(62, 3).Value = 0.7956586498444824
